# Update preferences order evaluation values (Mean/Std columns H/I)
# for PreOrder / PartialOrder preference-order metric rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Range("H26").Value = 0.77404
$ws.Range("I26").Value = 0.00754
$ws.Range("H27").Value = 0.1661
$ws.Range("I27").Value = 0.019
$ws.Range("H28").Value = 0.78443
$ws.Range("I28").Value = 0.00783
$ws.Range("H29").Value = 0.27578
$ws.Range("I29").Value = 0.02574
$ws.Range("H30").Value = 0.77413
$ws.Range("I30").Value = 0.00717
$ws.Range("H31").Value = 0.16502
$ws.Range("I31").Value = 0.01925
$ws.Range("H32").Value = 0.78671
$ws.Range("I32").Value = 0.008279999999999999
$ws.Range("H33").Value = 0.28259
$ws.Range("I33").Value = 0.02585
$ws.Range("H34").Value = 0.71472
$ws.Range("I34").Value = 0.00733
$ws.Range("H35").Value = 0.02751
$ws.Range("I35").Value = 0.01332
$ws.Range("H36").Value = 0.71485
$ws.Range("I36").Value = 0.00727
$ws.Range("H37").Value = 0.03133
$ws.Range("I37").Value = 0.01373
$ws.Range("H38").Value = 0.71469
$ws.Range("I38").Value = 0.00731
$ws.Range("H39").Value = 0.02751
$ws.Range("I39").Value = 0.01332
$ws.Range("H40").Value = 0.71469
$ws.Range("I40").Value = 0.00725
$ws.Range("H41").Value = 0.03091
$ws.Range("I41").Value = 0.01349
$ws.Range("H66").Value = 0.70928
$ws.Range("I66").Value = 0.00814
$ws.Range("H67").Value = 0.06107
$ws.Range("I67").Value = 0.01261
$ws.Range("H68").Value = 0.72934
$ws.Range("I68").Value = 0.00774
$ws.Range("H69").Value = 0.16186
$ws.Range("I69").Value = 0.01616
$ws.Range("H70").Value = 0.7099299999999999
$ws.Range("I70").Value = 0.00804
$ws.Range("H71").Value = 0.06116
$ws.Range("I71").Value = 0.01266
$ws.Range("H72").Value = 0.73424
$ws.Range("I72").Value = 0.007990000000000001
$ws.Range("H73").Value = 0.1769
$ws.Range("I73").Value = 0.0177
$ws.Range("H74").Value = 0.69664
$ws.Range("I74").Value = 0.00733
$ws.Range("H75").Value = 0.01811
$ws.Range("I75").Value = 0.00763
$ws.Range("H76").Value = 0.69882
$ws.Range("I76").Value = 0.00726
$ws.Range("H77").Value = 0.02883
$ws.Range("I77").Value = 0.009690000000000001
$ws.Range("H78").Value = 0.69668
$ws.Range("I78").Value = 0.00728
$ws.Range("H79").Value = 0.01811
$ws.Range("I79").Value = 0.00763
$ws.Range("H80").Value = 0.6986599999999999
$ws.Range("I80").Value = 0.0072
$ws.Range("H81").Value = 0.02842
$ws.Range("I81").Value = 0.00959
$ws.Range("H106").Value = 0.77214
$ws.Range("I106").Value = 0.009339999999999999
$ws.Range("H107").Value = 0.16336
$ws.Range("I107").Value = 0.02549
$ws.Range("H108").Value = 0.78126
$ws.Range("I108").Value = 0.00928
$ws.Range("H109").Value = 0.25784
$ws.Range("I109").Value = 0.03425
$ws.Range("H110").Value = 0.7723
$ws.Range("I110").Value = 0.009469999999999999
$ws.Range("H111").Value = 0.16262
$ws.Range("I111").Value = 0.02547
$ws.Range("H112").Value = 0.78327
$ws.Range("I112").Value = 0.00962
$ws.Range("H113").Value = 0.26416
$ws.Range("I113").Value = 0.03447
$ws.Range("H114").Value = 0.71283
$ws.Range("I114").Value = 0.00751
$ws.Range("H115").Value = 0.02459
$ws.Range("I115").Value = 0.01128
$ws.Range("H116").Value = 0.71306
$ws.Range("I116").Value = 0.00761
$ws.Range("H117").Value = 0.02684
$ws.Range("I117").Value = 0.01145
$ws.Range("H118").Value = 0.71284
$ws.Range("I118").Value = 0.00745
$ws.Range("H119").Value = 0.02459
$ws.Range("I119").Value = 0.01128
$ws.Range("H120").Value = 0.7129799999999999
$ws.Range("I120").Value = 0.00757
$ws.Range("H121").Value = 0.02667
$ws.Range("I121").Value = 0.01139
$ws.Range("H146").Value = 0.70969
$ws.Range("I146").Value = 0.00912
$ws.Range("H147").Value = 0.06647
$ws.Range("I147").Value = 0.01102
$ws.Range("H148").Value = 0.72396
$ws.Range("I148").Value = 0.0086
$ws.Range("H149").Value = 0.16028
$ws.Range("I149").Value = 0.01982
$ws.Range("H150").Value = 0.70958
$ws.Range("I150").Value = 0.00925
$ws.Range("H151").Value = 0.06647
$ws.Range("I151").Value = 0.01102
$ws.Range("H152").Value = 0.73064
$ws.Range("I152").Value = 0.00929
$ws.Range("H153").Value = 0.18164
$ws.Range("I153").Value = 0.02435
$ws.Range("H154").Value = 0.69723
$ws.Range("I154").Value = 0.008070000000000001
$ws.Range("H155").Value = 0.02086
$ws.Range("I155").Value = 0.00716
$ws.Range("H156").Value = 0.6986599999999999
$ws.Range("I156").Value = 0.008109999999999999
$ws.Range("H157").Value = 0.03191
$ws.Range("I157").Value = 0.00797
$ws.Range("H158").Value = 0.69722
$ws.Range("I158").Value = 0.008019999999999999
$ws.Range("H159").Value = 0.02086
$ws.Range("I159").Value = 0.00716
$ws.Range("H160").Value = 0.69856
$ws.Range("I160").Value = 0.00793
$ws.Range("H161").Value = 0.03157
$ws.Range("I161").Value = 0.00766
